$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.002957
$ws.Range("H2").Value = 60.008871
$ws.Range("I2").Value = 0.7920860939997775
$ws.Range("J2").Value = 0.7920860939997775
$ws.Range("M2").Value = 61.156892
$ws.Range("N2").Value = 183.470676
$ws.Range("O2").Value = 0.9308124486389074
$ws.Range("P2").Value = 0.9308124486389074
$ws.Range("Q2").Value = 1223.318680929644
$ws.Range("R2").Value = 11009.8681283668
$ws.Range("S2").Value = 0.7372835966887606
$ws.Range("T2").Value = 0.7372835966887606
$ws.Range("G3").Value = 20.002957
$ws.Range("H3").Value = 60.008871
$ws.Range("I3").Value = 0.7920860939997775
$ws.Range("J3").Value = 0.7920860939997775
$ws.Range("O3").Value = 0.02171808228502914
$ws.Range("P3").Value = 0.02171808228502914
$ws.Range("Q3").Value = 28.54295278505666
$ws.Range("R3").Value = 256.88657506551
$ws.Range("S3").Value = 0.01720259096631449
$ws.Range("T3").Value = 0.01720259096631449
$ws.Range("G4").Value = 20.002957
$ws.Range("H4").Value = 60.008871
$ws.Range("I4").Value = 0.7920860939997775
$ws.Range("J4").Value = 0.7920860939997775
$ws.Range("M4").Value = 2.00294
$ws.Range("N4").Value = 6.00882
$ws.Range("O4").Value = 0.03048489589491914
$ws.Range("P4").Value = 0.03048489589491914
$ws.Range("Q4").Value = 40.06472269358
$ws.Range("R4").Value = 360.58250424222
$ws.Range("S4").Value = 0.02414666211539635
$ws.Range("T4").Value = 0.02414666211539635
$ws.Range("G5").Value = 20.002957
$ws.Range("H5").Value = 60.008871
$ws.Range("I5").Value = 0.7920860939997775
$ws.Range("J5").Value = 0.7920860939997775
$ws.Range("M5").Value = 1.115932333333334
$ws.Range("N5").Value = 3.347797
$ws.Range("O5").Value = 0.01698457318114416
$ws.Range("P5").Value = 0.01698457318114415
$ws.Range("Q5").Value = 22.32194647857634
$ws.Range("R5").Value = 200.897518307187
$ws.Range("S5").Value = 0.01345324422930585
$ws.Range("T5").Value = 0.01345324422930585
$ws.Range("I6").Value = 0.04149178396178559
$ws.Range("J6").Value = 0.04149178396178559
$ws.Range("M6").Value = 61.156892
$ws.Range("N6").Value = 183.470676
$ws.Range("O6").Value = 0.9308124486389074
$ws.Range("P6").Value = 0.9308124486389074
$ws.Range("Q6").Value = 64.08100686282665
$ws.Range("R6").Value = 576.72906176544
$ws.Range("S6").Value = 0.0386210690278662
$ws.Range("T6").Value = 0.0386210690278662
$ws.Range("I7").Value = 0.04149178396178559
$ws.Range("J7").Value = 0.04149178396178559
$ws.Range("O7").Value = 0.02171808228502914
$ws.Range("P7").Value = 0.02171808228502914
$ws.Range("S7").Value = 0.0009011219782347119
$ws.Range("T7").Value = 0.0009011219782347118
$ws.Range("I8").Value = 0.04149178396178559
$ws.Range("J8").Value = 0.04149178396178559
$ws.Range("M8").Value = 2.00294
$ws.Range("N8").Value = 6.00882
$ws.Range("O8").Value = 0.03048489589491914
$ws.Range("P8").Value = 0.03048489589491914
$ws.Range("Q8").Value = 2.098707237866667
$ws.Range("R8").Value = 18.8883651408
$ws.Range("S8").Value = 0.00126487271456951
$ws.Range("T8").Value = 0.001264872714569509
$ws.Range("I9").Value = 0.04149178396178559
$ws.Range("J9").Value = 0.04149178396178559
$ws.Range("M9").Value = 1.115932333333334
$ws.Range("N9").Value = 3.347797
$ws.Range("O9").Value = 0.01698457318114416
$ws.Range("P9").Value = 0.01698457318114415
$ws.Range("Q9").Value = 1.169288777964445
$ws.Range("R9").Value = 10.52359900168
$ws.Range("S9").Value = 0.0007047202411151708
$ws.Range("T9").Value = 0.0007047202411151707
$ws.Range("G10").Value = 2.298356333333333
$ws.Range("H10").Value = 6.895068999999999
$ws.Range("I10").Value = 0.09101134850660582
$ws.Range("J10").Value = 0.09101134850660582
$ws.Range("M10").Value = 61.156892
$ws.Range("N10").Value = 183.470676
$ws.Range("O10").Value = 0.9308124486389074
$ws.Range("P10").Value = 0.9308124486389074
$ws.Range("Q10").Value = 140.5603300551826
$ws.Range("R10").Value = 1265.042970496644
$ws.Range("S10").Value = 0.08471449615736273
$ws.Range("T10").Value = 0.08471449615736273
$ws.Range("G11").Value = 2.298356333333333
$ws.Range("H11").Value = 6.895068999999999
$ws.Range("I11").Value = 0.09101134850660582
$ws.Range("J11").Value = 0.09101134850660582
$ws.Range("O11").Value = 0.02171808228502914
$ws.Range("P11").Value = 0.02171808228502914
$ws.Range("Q11").Value = 3.279608925098888
$ws.Range("R11").Value = 29.51648032589
$ws.Range("S11").Value = 0.001976591955737929
$ws.Range("T11").Value = 0.001976591955737929
$ws.Range("G12").Value = 2.298356333333333
$ws.Range("H12").Value = 6.895068999999999
$ws.Range("I12").Value = 0.09101134850660582
$ws.Range("J12").Value = 0.09101134850660582
$ws.Range("M12").Value = 2.00294
$ws.Range("N12").Value = 6.00882
$ws.Range("O12").Value = 0.03048489589491914
$ws.Range("P12").Value = 0.03048489589491914
$ws.Range("Q12").Value = 4.603469834286666
$ws.Range("R12").Value = 41.43122850858
$ws.Range("S12").Value = 0.002774471484480083
$ws.Range("T12").Value = 0.002774471484480083
$ws.Range("G13").Value = 2.298356333333333
$ws.Range("H13").Value = 6.895068999999999
$ws.Range("I13").Value = 0.09101134850660582
$ws.Range("J13").Value = 0.09101134850660582
$ws.Range("M13").Value = 1.115932333333334
$ws.Range("N13").Value = 3.347797
$ws.Range("O13").Value = 0.01698457318114416
$ws.Range("P13").Value = 0.01698457318114415
$ws.Range("Q13").Value = 2.564810145888111
$ws.Range("R13").Value = 23.083291312993
$ws.Range("S13").Value = 0.001545788909025062
$ws.Range("T13").Value = 0.001545788909025061
$ws.Range("G14").Value = 1.904387
$ws.Range("H14").Value = 5.713160999999999
$ws.Range("I14").Value = 0.07541077353183102
$ws.Range("J14").Value = 0.07541077353183102
$ws.Range("M14").Value = 61.156892
$ws.Range("N14").Value = 183.470676
$ws.Range("O14").Value = 0.9308124486389074
$ws.Range("P14").Value = 0.9308124486389074
$ws.Range("Q14").Value = 116.466390085204
$ws.Range("R14").Value = 1048.197510766836
$ws.Range("S14").Value = 0.07019328676491775
$ws.Range("T14").Value = 0.07019328676491775
$ws.Range("G15").Value = 1.904387
$ws.Range("H15").Value = 5.713160999999999
$ws.Range("I15").Value = 0.07541077353183102
$ws.Range("J15").Value = 0.07541077353183102
$ws.Range("O15").Value = 0.02171808228502914
$ws.Range("P15").Value = 0.02171808228502914
$ws.Range("Q15").Value = 2.717439637823333
$ws.Range("R15").Value = 24.45695674041
$ws.Range("S15").Value = 0.001637777384742004
$ws.Range("T15").Value = 0.001637777384742004
$ws.Range("G16").Value = 1.904387
$ws.Range("H16").Value = 5.713160999999999
$ws.Range("I16").Value = 0.07541077353183102
$ws.Range("J16").Value = 0.07541077353183102
$ws.Range("M16").Value = 2.00294
$ws.Range("N16").Value = 6.00882
$ws.Range("O16").Value = 0.03048489589491914
$ws.Range("P16").Value = 0.03048489589491914
$ws.Range("Q16").Value = 3.81437289778
$ws.Range("R16").Value = 34.32935608002
$ws.Range("S16").Value = 0.002298889580473193
$ws.Range("T16").Value = 0.002298889580473193
$ws.Range("G17").Value = 1.904387
$ws.Range("H17").Value = 5.713160999999999
$ws.Range("I17").Value = 0.07541077353183102
$ws.Range("J17").Value = 0.07541077353183102
$ws.Range("M17").Value = 1.115932333333334
$ws.Range("N17").Value = 3.347797
$ws.Range("O17").Value = 0.01698457318114416
$ws.Range("P17").Value = 0.01698457318114415
$ws.Range("Q17").Value = 2.125167028479667
$ws.Range("R17").Value = 19.126503256317
$ws.Range("S17").Value = 0.001280819801698073
$ws.Range("T17").Value = 0.001280819801698073
